$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.951.24'
$ws.Range("E2").Value = '  +4.52%  '
$ws.Range("D3").Value = '3.240.13'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'395.09"
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").Value = "'108.13"
$ws.Range("E6").Value = '  -2.10%  '
$ws.Range("D7").Value = "'0.582"
$ws.Range("E7").Value = '  +6.06%  '
$ws.Range("D8").Value = '3.235.43'
$ws.Range("E8").Value = '  +2.03%  '
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = "'0.627"
$ws.Range("E10").Value = '  +1.29%  '
$ws.Range("D11").Value = "'39.22"
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = "'0.0971"
$ws.Range("E12").Value = '  +9.00%  '
$ws.Range("D14").Value = '3.751.69'
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("D15").Value = "'8.41"
$ws.Range("E15").Value = '  +4.42%  '
$ws.Range("D16").Value = "'18.96"
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '3.241.65'
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("E18").Value = '  -3.08%  '
$ws.Range("D19").Value = "'10.99"
$ws.Range("E19").Value = '  +4.17%  '
$ws.Range("D20").Value = '56.836.31'
$ws.Range("E20").Value = '  +4.44%  '
$ws.Range("D21").Value = "'3.31"
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").Value = "'0.0000110"
$ws.Range("E22").Value = '  +9.88%  '
$ws.Range("D23").Value = "'13.09"
$ws.Range("E23").Value = '  +1.68%  '
$ws.Range("D24").Value = "'290.56"
$ws.Range("E24").Value = '  +5.60%  '
$ws.Range("D25").Value = "'73.85"
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("E26").Value = '  -3.54%  '
$ws.Range("D27").Value = "'27.97"
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = "'4.35"
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("D29").Value = "'7.69"
$ws.Range("E29").Value = '  -4.80%  '
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").Value = "'7.14"
$ws.Range("E32").Value = '  -5.29%  '
$ws.Range("D33").Value = "'11.17"
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("E34").Value = '  -2.47%  '
$ws.Range("D35").Value = "'39.98"
$ws.Range("E35").Value = '  +9.03%  '
$ws.Range("D36").Value = "'0.0482"
$ws.Range("E36").Value = '  -2.73%  '
$ws.Range("E37").Value = '  +1.49%  '
$ws.Range("D38").Value = "'51.43"
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("E40").Value = '  -4.05%  '
$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("D42").Value = "'137.31"
$ws.Range("E42").Value = '  +3.41%  '
$ws.Range("D43").Value = "'0.122"
$ws.Range("E43").Value = '  +4.26%  '
$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").Value = "'4.03"
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").Value = "'1.89"
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("D46").Value = "'17.01"
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = "'0.274"
$ws.Range("E47").Value = '  -4.94%  '
$ws.Range("D48").Value = "'2.30"
$ws.Range("E48").Value = '  +10.49%  '
$ws.Range("D49").Value = "'22.22"
$ws.Range("E49").Value = '  +0.35%  '
$ws.Range("D50").Value = '2.166.49'
$ws.Range("E50").Value = '  +3.12%  '
$ws.Range("E51").Value = '  -5.62%  '
